$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1540.6786
$ws.Range("J17").Value = 1540.6786
$ws.Range("L17").Value = 4622.0358
$ws.Range("N17").Value = -4958.0358

$ws.Range("H100").Value = 963.3333
$ws.Range("I100").Value = 945.75
$ws.Range("J100").Value = 998.5
$ws.Range("K100").Value = 945.75
$ws.Range("L100").Value = 998.5
$ws.Range("M100").Value = -404.75
$ws.Range("N100").Value = -2080.5

$ws.Range("H111").Value = 343
$ws.Range("I111").Value = 364.5
$ws.Range("K111").Value = 1093.5
$ws.Range("M111").Value = 1973.5

$ws.Range("H113").Value = 5683.273
$ws.Range("I113").Value = 4740.2856
$ws.Range("J113").Value = 7333.5
$ws.Range("K113").Value = 4740.2856
$ws.Range("L113").Value = 7333.5
$ws.Range("M113").Value = -1486.2856
$ws.Range("N113").Value = -13841.5

$ws.Range("H138").Value = 2767.494
$ws.Range("I138").Value = 3584.182
$ws.Range("J138").Value = 2462.966
$ws.Range("K138").Value = 10752.546
$ws.Range("L138").Value = 7388.897999999999
$ws.Range("M138").Value = -5612.545999999998
$ws.Range("N138").Value = -17668.898

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1692.4546
$ws.Range("I2").Value = 1711.7
$ws.Range("K2").Value = 1711.7
$ws.Range("M2").Value = -1598.7

$ws.Range("H97").Value = 1012.36365
$ws.Range("I97").Value = 584.381
$ws.Range("K97").Value = 584.381
$ws.Range("M97").Value = -88.38099999999997

$ws.Range("H102").Value = 1004
$ws.Range("I102").Value = 1010
$ws.Range("J102").Value = 998
$ws.Range("K102").Value = 1010
$ws.Range("L102").Value = 998
$ws.Range("M102").Value = 612
$ws.Range("N102").Value = -4242

$ws.Range("H116").Value = 1692.4546
$ws.Range("I116").Value = 1711.7
$ws.Range("K116").Value = 1711.7
$ws.Range("M116").Value = 582.3

$ws.Range("H140").Value = 67500
$ws.Range("J140").Value = 67500
$ws.Range("L140").Value = 67500
$ws.Range("N140").Value = -77860

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1692.4546
$ws.Range("I3").Value = 1711.7
$ws.Range("K3").Value = 1711.7
$ws.Range("M3").Value = -1597.7

$ws.Range("H54").Value = 15000
$ws.Range("I54").Value = 15000
$ws.Range("K54").Value = 15000
$ws.Range("M54").Value = -14516

$ws.Range("H94").Value = 1443.6316
$ws.Range("I94").Value = 732.64703
$ws.Range("J94").Value = 7487
$ws.Range("K94").Value = 732.64703
$ws.Range("L94").Value = 7487
$ws.Range("M94").Value = -281.64703
$ws.Range("N94").Value = -8389

$ws.Range("H107").Value = 905.7143
$ws.Range("I107").Value = 821.53845
$ws.Range("K107").Value = 821.53845
$ws.Range("M107").Value = 1098.46155

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1622
$ws.Range("I16").Value = 744
$ws.Range("K16").Value = 744
$ws.Range("M16").Value = -457

$ws.Range("H105").Value = 2908.7144
$ws.Range("I105").Value = 2740.5
$ws.Range("K105").Value = 2740.5
$ws.Range("M105").Value = -993.5

$ws.Range("H107").Value = 800
$ws.Range("I107").Value = 800
$ws.Range("K107").Value = 800
$ws.Range("M107").Value = 1120

$ws.Range("H113").Value = 1622
$ws.Range("I113").Value = 744
$ws.Range("K113").Value = 744
$ws.Range("M113").Value = 1426

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 333372.72
$ws.Range("I2").Value = 454577.47
$ws.Range("J2").Value = 59.75
$ws.Range("K2").Value = 2727464.82
$ws.Range("L2").Value = 358.5
$ws.Range("M2").Value = -2727351.82
$ws.Range("N2").Value = -584.5

$ws.Range("H116").Value = 25987.5
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1513.75
$ws.Range("J97").Value = 2898
$ws.Range("L97").Value = 2898
$ws.Range("N97").Value = -3890

$ws.Range("H113").Value = 3392.9092
$ws.Range("I113").Value = 2517.5715
$ws.Range("J113").Value = 4924.75
$ws.Range("K113").Value = 2517.5715
$ws.Range("L113").Value = 4924.75
$ws.Range("M113").Value = -347.5715
$ws.Range("N113").Value = -9264.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1460.7142
$ws.Range("I55").Value = 333.1579
$ws.Range("K55").Value = 333.1579
$ws.Range("M55").Value = -160.1579

$ws.Range("H93").Value = 2526.9285
$ws.Range("I93").Value = 2526.9285
$ws.Range("K93").Value = 2526.9285
$ws.Range("M93").Value = -1278.9285

$ws.Range("H109").Value = 82999.5
$ws.Range("J109").Value = 82999.5
$ws.Range("L109").Value = 82999.5
$ws.Range("N109").Value = -85773.5

$ws.Range("H119").Value = 64613.332
$ws.Range("J119").Value = 64613.332
$ws.Range("L119").Value = 64613.332
$ws.Range("N119").Value = -74289.33199999999

$ws.Range("H120").Value = 89000
$ws.Range("J120").Value = 89000
$ws.Range("L120").Value = 89000
$ws.Range("N120").Value = -98676

$ws.Range("H121").Value = 79999.336
$ws.Range("J121").Value = 79999.336
$ws.Range("L121").Value = 79999.336
$ws.Range("N121").Value = -83493.336

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H128").Value = 70482.25
$ws.Range("J128").Value = 70482.25
$ws.Range("L128").Value = 70482.25
$ws.Range("N128").Value = -80442.25

$ws.Range("H130").Value = 63250
$ws.Range("J130").Value = 63250
$ws.Range("L130").Value = 63250
$ws.Range("N130").Value = -73290

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 480.2857
$ws.Range("I100").Value = 460.33334
$ws.Range("K100").Value = 920.66668
$ws.Range("M100").Value = -379.66668

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws.Range("H123").Value = 80250
$ws.Range("J123").Value = 80250
$ws.Range("L123").Value = 80250
$ws.Range("N123").Value = -90050

$ws.Range("H126").Value = 966.6667
$ws.Range("I126").Value = 980
$ws.Range("J126").Value = 900
$ws.Range("K126").Value = 2940
$ws.Range("L126").Value = 2700
$ws.Range("M126").Value = -470
$ws.Range("N126").Value = -7640
